# "Generate Report for Handback"
#
# The localization-status report is regenerated after the de-de / zh-cn
# handback files come back in sync with en-US:
#   - Overview sheet "Status" column flips from "Ready for handoff" to
#     "Handed back: in sync with en-US".
#   - Each language sheet's "Latest Target File" / "Latest Handback File"
#     columns (I/J) get populated (they were blank while waiting on
#     handback), with a hyperlink on the target-file cell mirroring the
#     source-file hyperlink in column A.
#   - de-de's "Latest Handback DateTime" column (K) gets a real timestamp
#     now that the handback actually happened.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$mdFile1 = "4674edd2-9dff-47e4-b740-37ac1cfed0ee.md"
$mdFile2 = "83ff810f-bfd3-4f3e-80bb-0e5686272718.md"
$mdUrl1  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/297f0147d86be730960dde2a9293512168b32ec7/e2e/$mdFile1"
$mdUrl2  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/297f0147d86be730960dde2a9293512168b32ec7/e2e/$mdFile2"

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: status column (E/F) for both rows, plus column autofit
# ---------------------------------------------------------------------
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527

# ---------------------------------------------------------------------
# zh-cn sheet: target/handback file columns now populated
# ---------------------------------------------------------------------
$zhcn.Range("I2").Value = $mdFile1
$zhcn.Range("I2").Style = "HyperLink"
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $mdUrl1, [System.Type]::Missing, [System.Type]::Missing, $mdFile1) | Out-Null
$zhcn.Range("J2").Value = "4674edd2-9dff-47e4-b740-37ac1cfed0ee.257924c51840b859ed55eae9f49f978c450d833f.zh-cn.xlf"

$zhcn.Range("I3").Value = $mdFile2
$zhcn.Range("I3").Style = "HyperLink"
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $mdUrl2, [System.Type]::Missing, [System.Type]::Missing, $mdFile2) | Out-Null
$zhcn.Range("J3").Value = "83ff810f-bfd3-4f3e-80bb-0e5686272718.f10e8c1dca65599f017eeb250abb1441f332aabf.zh-cn.xlf"

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$zhcn.Columns.Item(3).ColumnWidth = 29.9777047293527
# Excel stores ColumnWidth internally as whole pixels (MDW=6 + 5px
# padding), so asking for an even "40" lands on 40.8333; dialing the
# input back by the padding amount lands the stored width on exactly 40.
$zhcn.Columns.Item(9).ColumnWidth = 39.1666666666667
$zhcn.Columns.Item(10).ColumnWidth = 39.1666666666667

# ---------------------------------------------------------------------
# de-de sheet: target/handback file columns now populated, and the
# handback actually completed so the datetime moves off the epoch
# sentinel.
# ---------------------------------------------------------------------
$dede.Range("I2").Value = $mdFile1
$dede.Range("I2").Style = "HyperLink"
$dede.Hyperlinks.Add($dede.Range("I2"), $mdUrl1, [System.Type]::Missing, [System.Type]::Missing, $mdFile1) | Out-Null
$dede.Range("J2").Value = "4674edd2-9dff-47e4-b740-37ac1cfed0ee.257924c51840b859ed55eae9f49f978c450d833f.de-de.xlf"
$dede.Range("K2").Value = "2016-09-07 04:40:00"

$dede.Range("I3").Value = $mdFile2
$dede.Range("I3").Style = "HyperLink"
$dede.Hyperlinks.Add($dede.Range("I3"), $mdUrl2, [System.Type]::Missing, [System.Type]::Missing, $mdFile2) | Out-Null
$dede.Range("J3").Value = "83ff810f-bfd3-4f3e-80bb-0e5686272718.f10e8c1dca65599f017eeb250abb1441f332aabf.de-de.xlf"
$dede.Range("K3").Value = "2016-09-07 04:40:00"

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

$dede.Columns.Item(3).ColumnWidth = 29.9777047293527
$dede.Columns.Item(9).ColumnWidth = 39.1666666666667
$dede.Columns.Item(10).ColumnWidth = 39.1666666666667

Write-Output "Handback report regenerated."
